$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set new/updated values in column F (săpt. 4) for several students
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 2

# Update the active selection to match the author's final view state
$ws.Range("I6").Select()
